$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "ImageName"
$ws.Range("H2").Value = "exm.jpg"

$ws.Range("H2").Select()
